$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text (it stores values like "27.455.31" which
# are not valid numbers, alongside plain-decimal-looking values like "1.005"
# that Excel would otherwise auto-convert to a number). Apply a text number
# format before writing, then restore the default "Normal" style afterward so
# no stray per-cell formatting is left behind.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.450.52'
$ws.Range("E2").Value = '  -1.27%  '
$ws.Range("D3").Value = '1.734.87'
$ws.Range("E3").Value = '  -1.47%  '
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.52%  '
$ws.Range("D5").Value = '321.95'
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").Value = '  +0.51%  '
$ws.Range("D7").Value = '0.4620'
$ws.Range("E7").Value = '  +8.95%  '
$ws.Range("D8").Value = '0.3521'
$ws.Range("E8").Value = '  -2.90%  '
$ws.Range("D9").Value = '41.70'
$ws.Range("E9").Value = '  -1.81%  '
$ws.Range("D10").Value = '0.07350'
$ws.Range("E10").Value = '  -1.61%  '
$ws.Range("D11").Value = '1.076'
$ws.Range("E11").Value = '  -0.95%  '
$ws.Range("D12").Value = '1.004'
$ws.Range("E12").Value = '  +0.46%  '
$ws.Range("D13").Value = '20.42'
$ws.Range("E13").Value = '  -1.19%  '
$ws.Range("D14").Value = '5.910'
$ws.Range("E14").Value = '  -2.57%  '
$ws.Range("D15").Value = '7.043'
$ws.Range("E15").Value = '  -3.41%  '
$ws.Range("D16").Value = '1.736.44'
$ws.Range("E16").Value = '  -1.06%  '
$ws.Range("D17").Value = '91.00'
$ws.Range("E17").Value = '  +0.04%  '
$ws.Range("D18").Value = '0.00001052'
$ws.Range("E18").Value = '  -0.26%  '
$ws.Range("D19").Value = '0.06390'
$ws.Range("E19").Value = '  +0.32%  '
$ws.Range("D20").Value = '1.004'
$ws.Range("E20").Value = '  +0.42%  '
$ws.Range("D21").Value = '16.62'
$ws.Range("E21").Value = '  -2.24%  '
$ws.Range("D22").Value = '5.723'
$ws.Range("D23").Value = '27.511.36'
$ws.Range("E23").Value = '  -1.11%  '
$ws.Range("E24").Value = '  -1.20%  '
$ws.Range("D25").Value = '2.095'
$ws.Range("E25").Value = '  -0.23%  '
$ws.Range("D26").Value = '162.59'
$ws.Range("E26").Value = '  +3.28%  '
$ws.Range("E27").Value = '  -1.83%  '
$ws.Range("D28").Value = '1.937.49'
$ws.Range("E28").Value = '  -1.08%  '
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").Value = '2.035'
$ws.Range("E29").Value = '  -4.57%  '
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").Value = '124.34'
$ws.Range("E30").Value = '  +0.23%  '
$ws.Range("E31").Value = '  -6.41%  '
$ws.Range("D32").Value = '0.09193'
$ws.Range("E32").Value = '  +3.88%  '
$ws.Range("D33").Value = '3.667'
$ws.Range("E33").Value = '  -0.39%  '
$ws.Range("D34").Value = '5.404'
$ws.Range("E34").Value = '  -2.75%  '
$ws.Range("D35").Value = '0.02265'
$ws.Range("E35").Value = '  -0.99%  '
$ws.Range("D36").Value = '11.58'
$ws.Range("E36").Value = '  -5.23%  '
$ws.Range("D37").Value = '0.05981'
$ws.Range("E37").Value = '  -1.06%  '
$ws.Range("D38").Value = '0.2064'
$ws.Range("E38").Value = '  -1.53%  '
$ws.Range("D39").Value = '4.907'
$ws.Range("E39").Value = '  -1.17%  '
$ws.Range("D40").Value = '0.6229'
$ws.Range("E40").Value = '  -1.40%  '
$ws.Range("D41").Value = '1.184'
$ws.Range("E41").Value = '  +0.77%  '
$ws.Range("D42").Value = '1.375'
$ws.Range("E42").Value = '  -1.47%  '
$ws.Range("D43").Value = '7.713'
$ws.Range("E43").Value = '  -2.02%  '
$ws.Range("D44").Value = '13.04'
$ws.Range("E44").Value = '  -1.47%  '
$ws.Range("D45").Value = '3.696'
$ws.Range("E45").Value = '  +0.38%  '
$ws.Range("D46").Value = '0.5792'
$ws.Range("E46").Value = '  -1.19%  '
$ws.Range("D47").Value = '121.74'
$ws.Range("E47").Value = '  -0.93%  '
$ws.Range("E48").Value = '  -3.23%  '
$ws.Range("D49").Value = '0.06828'
$ws.Range("E49").Value = '  +0.11%  '
$ws.Range("D50").Value = '1.121'
$ws.Range("E50").Value = '  -4.90%  '
$ws.Range("D51").Value = '71.14'
$ws.Range("E51").Value = '  -3.41%  '

# Restore default styling on column D so the cells match the original
# (un-styled) appearance; only the underlying values should have changed.
$ws.Range("D2:D51").Style = "Normal"
